$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: the heading run " Control & Computing System Design:" that
# follows "Safety Steps" becomes " in Control & Computing System
# Design:" so the heading reads
# "Safety Steps in Control & Computing System Design:"
# ------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute(" Control & Computing System Design:", $false, $false, $false, $false, $false, $true, 1, $false, " in Control & Computing System Design:", 2)

# ------------------------------------------------------------------
# Change 2: add a new "Safety Steps in Propulsion Design:" heading and
# its body paragraph right after the Control & Computing System
# Design section (i.e. right after the paragraph that ends with
# "...without any glitch or unwanted output.")
# ------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("without any glitch or unwanted output.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lastPara = $find2.Paragraphs(1)

# Insert a new (initially empty) paragraph right after it and turn it
# into the "Safety Steps in Propulsion Design:" heading.
$lastPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs(6)
$headingPara.Range.Text = "Safety Steps in Propulsion Design:"
$headingPara.Style = "Heading1"

# Insert another new (initially empty) paragraph right after the new
# heading and fill it in with the Propulsion write-up, using raw OOXML
# so the leading tab becomes a real <w:tab/> run (matching how the
# rest of the document represents its paragraph-leading tabs).
$headingPara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs(7)
$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Safer propulsion system implies safer flight. Motors &amp; Propellers used in the drone are the highest quality ones available in market to ensure the safest propulsion system possible.</w:t></w:r></w:p>'
$bodyPara.Range.InsertXML($bodyXml)
